# "Excel file updated after additions of Post02"
#
# The "Log of all Blogs" table (Table2, on Sheet1) gets one new row added
# for the "Master CSS :- Top CSS Frameworks and Resources" blog post:
#   B12 = 2 (S.No)
#   C12 = title
#   D12 = date of post (2020-10-09)
#   E12 = hyperlink to the Hashnode post
#   F12 = hyperlink to the Dev.to post
# The table range / autofilter / dimension grow from B10:F11 to B10:F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table2")

# Grow the table by one row - this keeps the ListObject (and therefore the
# worksheet dimension / autofilter ref) in sync automatically.
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

# Carry the formatting (number format, hyperlink style, etc.) down from the
# row above, same as Excel does when you extend a table, then overwrite the
# values for the new row.
$ws.Range("B11:F11").Copy() | Out-Null
$rng.PasteSpecial(-4122) | Out-Null

$rng.Cells.Item(1, 1).Value = 2
$rng.Cells.Item(1, 2).Value = "Master CSS :- Top CSS Frameworks and Resources"
$rng.Cells.Item(1, 3).Value = "10/9/2020"

$ws.Hyperlinks.Add($rng.Cells.Item(1, 4), "https://programmingport.hashnode.dev/master-css-top-css-frameworks-and-resources", "", "", "") | Out-Null
$ws.Hyperlinks.Add($rng.Cells.Item(1, 5), "https://dev.to/rahulmishra05/master-css-top-css-frameworks-and-resources-5gj9", "", "", "") | Out-Null
